$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")
$ws.Activate()

# Update the end-time on row 9 (Tidspunkt til) from 16:00 to 16:45 -
# a 45 minute increase (0.75h / 24h = 0.03125 of a day).
$ws.Range("E9").Value = 16/24 + 45/1440

# Update the active selection to match the author's final cursor position.
$ws.Range("C9").Select()
